$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update STATUS (column B) values
$ws.Range("B73").Value = "SUCESSO"
$ws.Range("B92").Value = "ERRO"

# Update DATA EXTRACAO (column C) values
$ws.Range("C2").Value = "2025-02-02 20:00:31"
$ws.Range("C3").Value = "2025-02-02 20:00:41"
$ws.Range("C4").Value = "2025-02-02 20:01:29"
$ws.Range("C5").Value = "2025-02-02 20:01:34"
$ws.Range("C6").Value = "2025-02-02 20:01:47"
$ws.Range("C7").Value = "2025-02-02 20:02:12"
$ws.Range("C8").Value = "2025-02-02 20:02:17"
$ws.Range("C9").Value = "2025-02-02 20:02:25"
$ws.Range("C10").Value = "2025-02-02 20:02:31"
$ws.Range("C11").Value = "2025-02-02 20:02:41"
$ws.Range("C12").Value = "2025-02-02 20:02:50"
$ws.Range("C13").Value = "2025-02-03 00:02:47"
$ws.Range("C14").Value = "2025-02-03 00:02:54"
$ws.Range("C15").Value = "2025-02-03 00:03:00"
$ws.Range("C16").Value = "2025-02-03 00:03:05"
$ws.Range("C17").Value = "2025-02-03 00:03:12"
$ws.Range("C18").Value = "2025-02-03 00:03:24"
$ws.Range("C19").Value = "2025-02-03 00:03:32"
$ws.Range("C20").Value = "2025-02-03 00:03:38"
$ws.Range("C55").Value = "2025-02-02 22:30:23"
$ws.Range("C56").Value = "2025-02-02 22:30:28"
$ws.Range("C57").Value = "2025-02-02 22:30:37"
$ws.Range("C58").Value = "2025-02-02 22:30:45"
$ws.Range("C59").Value = "2025-02-02 22:30:53"
$ws.Range("C60").Value = "2025-02-02 22:30:58"
$ws.Range("C61").Value = "2025-02-02 22:31:11"
$ws.Range("C62").Value = "2025-02-02 22:31:17"
$ws.Range("C63").Value = "2025-02-02 22:31:34"
$ws.Range("C64").Value = "2025-02-02 22:31:44"
$ws.Range("C65").Value = "2025-02-02 22:31:58"
$ws.Range("C66").Value = "2025-02-02 22:32:09"
$ws.Range("C67").Value = "2025-02-02 22:32:17"
$ws.Range("C68").Value = "2025-02-02 22:32:34"
$ws.Range("C69").Value = "2025-02-02 22:32:47"
$ws.Range("C70").Value = "2025-02-02 22:32:54"
$ws.Range("C71").Value = "2025-02-02 22:33:03"
$ws.Range("C72").Value = "2025-02-02 22:33:16"
$ws.Range("C73").Value = "2025-02-02 22:33:22"
$ws.Range("C74").Value = "2025-02-02 22:33:44"
$ws.Range("C75").Value = "2025-02-02 22:33:53"
$ws.Range("C76").Value = "2025-02-02 22:34:25"
$ws.Range("C77").Value = "2025-02-02 22:34:34"
$ws.Range("C78").Value = "2025-02-02 22:34:40"
$ws.Range("C79").Value = "2025-02-02 22:34:48"
$ws.Range("C80").Value = "2025-02-02 22:34:55"
$ws.Range("C81").Value = "2025-02-02 22:35:04"
$ws.Range("C82").Value = "2025-02-02 22:35:11"
$ws.Range("C83").Value = "2025-02-02 22:35:15"
$ws.Range("C84").Value = "2025-02-02 22:35:28"
$ws.Range("C85").Value = "2025-02-02 22:35:36"
$ws.Range("C86").Value = "2025-02-02 22:35:41"
$ws.Range("C87").Value = "2025-02-02 22:36:52"
$ws.Range("C88").Value = "2025-02-02 22:37:01"
$ws.Range("C89").Value = "2025-02-02 22:37:07"
$ws.Range("C90").Value = "2025-02-02 22:37:16"
$ws.Range("C91").Value = "2025-02-02 22:37:24"
$ws.Range("C92").Value = "2025-02-02 22:37:31"
$ws.Range("C93").Value = "2025-02-02 22:37:55"

Write-Host "Updated $($($ws.UsedRange.Rows.Count)) rows"
